$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the cells below can
# be updated, then re-protect once the edits are in place.
$ws.Unprotect()

# Update confidential disclaimer text date from 2021-04-22 to 2021-04-23
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# Update Weight and Percent Change values
$ws.Range("D2").Value = 0.2516038891816318
$ws.Range("E2").Value = 0.01642760818928535

$ws.Range("D3").Value = 0.2472328959490401
$ws.Range("E3").Value = 0.01871580765908454

$ws.Range("D4").Value = 0.2484263306629009
$ws.Range("E4").Value = 0.01094963169420682

$ws.Range("D5").Value = 0.2527368842064273
$ws.Range("E5").Value = 0.0116448326055314

$ws.Range("E6").Value = 0.01442366897168634

# Restore protection on the sheet.
$ws.Protect()
